$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 352 (shifts the existing rows 352:379 down to 353:380,
# preserving formatting/styles the way Excel's Insert does).
$ws.Rows.Item(352).Insert()

# Populate the newly inserted row 352 with the new record.
$ws.Cells.Item(352, 1).Value = 5
$ws.Cells.Item(352, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(352, 3).Value = "Maule"
$ws.Cells.Item(352, 4).Value = 44783
$ws.Cells.Item(352, 5).Value = 7
$ws.Cells.Item(352, 6).Value = 100114013
$ws.Cells.Item(352, 7).Value = "Zanahoria"
$ws.Cells.Item(352, 8).Value = "Sin especificar"
$ws.Cells.Item(352, 9).Value = "Primera"
$ws.Cells.Item(352, 10).Value = 500
$ws.Cells.Item(352, 11).Value = 9000
$ws.Cells.Item(352, 12).Value = 9000
$ws.Cells.Item(352, 13).Value = 9000
$ws.Cells.Item(352, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(352, 15).Value = "Región de Ñuble"
$ws.Cells.Item(352, 16).Value = 450
$ws.Cells.Item(352, 17).Value = 20
$ws.Cells.Item(352, 18).Value = "Hortaliza"
